$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.0000", "28.145.16").
# Force text format before assignment so Excel does not coerce it to a Double
# (which would drop trailing zeros / thousands-style dots), then restore the
# default "Normal" style so no stray s="n" attribute is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.164.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.783.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4916'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2676'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06254'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.777.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.45'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07026'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6259'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.636'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '79.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.141.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9995'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007214'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.006.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.563'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.727'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.236'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.859'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '109.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.387'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.167'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08268'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.774'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04890'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.071'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.612'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6513'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9434'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.587'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.047'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.947'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01550'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.0000'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.78'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3986'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.180'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1201'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05433'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.017'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.297'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.26%  '
